$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.871.72'
$ws.Range("E2").Value = '  -3.81%  '
$ws.Range("D3").Value = '3.335.92'
$ws.Range("E3").Value = '  -1.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '574.67'
$ws.Range("E5").Value = '  -3.28%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '183.09'
$ws.Range("E6").Value = '  -4.95%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.599'
$ws.Range("E8").Value = '  -1.61%  '
$ws.Range("E9").Value = '  -3.86%  '
$ws.Range("E10").Value = '  -1.63%  '
$ws.Range("E11").Value = '  -4.13%  '
$ws.Range("D12").Value = '3.913.82'
$ws.Range("E12").Value = '  -1.25%  '
$ws.Range("E13").Value = '  -0.94%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.22'
$ws.Range("E14").Value = '  -5.35%  '
$ws.Range("D15").Value = '66.897.72'
$ws.Range("E15").Value = '  -3.82%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000168'
$ws.Range("E16").Value = '  -2.44%  '
$ws.Range("D17").Value = '3.308.99'
$ws.Range("E17").Value = '  -2.34%  '
$ws.Range("B18").Value = 'BitcoinCash'
$ws.Range("C18").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '435.71'
$ws.Range("E18").Value = '  -3.38%  '
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.68'
$ws.Range("E19").Value = '  -0.81%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.71'
$ws.Range("E20").Value = '  -2.33%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.67'
$ws.Range("E21").Value = '  -1.98%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '73.74'
$ws.Range("E22").Value = '  +0.34%  '
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.521'
$ws.Range("E24").Value = '  +0.36%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000119'
$ws.Range("E25").Value = '  -2.21%  '
$ws.Range("E26").Value = '  -0.69%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.10'
$ws.Range("E27").Value = '  -5.16%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("E28").Value = '  -0.20%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.97'
$ws.Range("E29").Value = '  -1.79%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '22.88'
$ws.Range("E30").Value = '  -1.68%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.34'
$ws.Range("E31").Value = '  -5.09%  '
$ws.Range("B32").Value = 'Aptos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.86'
$ws.Range("E32").Value = '  -2.50%  '
$ws.Range("B33").Value = 'USDe'
$ws.Range("C33").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  +0.06%  '
$ws.Range("B34").Value = 'Fetch.AI'
$ws.Range("C34").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.24'
$ws.Range("E34").Value = '  -5.05%  '
$ws.Range("E35").Value = '  -0.73%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '160.06'
$ws.Range("E36").Value = '  -2.85%  '
$ws.Range("B37").Value = 'EnergySwap'
$ws.Range("C37").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '27.37'
$ws.Range("E37").Value = '  +0.46%  '
$ws.Range("B38").Value = 'Stacks'
$ws.Range("C38").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.85'
$ws.Range("E38").Value = '  -4.37%  '
$ws.Range("D39").Value = '2.837.24'
$ws.Range("E39").Value = '  +3.26%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.793'
$ws.Range("E40").Value = '  -3.68%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.46'
$ws.Range("E41").Value = '  -3.14%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.25'
$ws.Range("E42").Value = '  -4.17%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0678'
$ws.Range("E43").Value = '  -1.82%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.27'
$ws.Range("E44").Value = '  -1.38%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '24.66'
$ws.Range("E45").Value = '  -3.37%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.37'
$ws.Range("E46").Value = '  -6.33%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '324.19'
$ws.Range("E47").Value = '  -5.99%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0274'
$ws.Range("E48").Value = '  -4.10%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.994'
$ws.Range("E49").Value = '  -3.91%  '
$ws.Range("E50").Value = '  -2.47%  '
$ws.Range("E51").Value = '  -1.38%  '